# Applies updated price/volume data for the symbol list refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.11%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'36.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.24%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.063"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.06%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08122"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.26%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.084"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'8.38%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.158"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.11%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'7.861"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.24%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9295"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.29%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1416"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'11.63%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.49%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09116"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.89%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03458"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.05%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09917"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.001407"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.79%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006320"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.04%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.839"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'6.20%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.334"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.23%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3447"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.12%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1284"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.47%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.805"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-7.14%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2341"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-7.44%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04363"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.86%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.21%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'3.90%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'0.0001298"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.12%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02018"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.04%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05170"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.81%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007482"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.84%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.34%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1370"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.20%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002167"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.21%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009980"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.67%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006267"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.94%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'64.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.56%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001251"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-21.80%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.09%"
$ws.Range("E51").Style = "Normal"
